$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) values
$ws.Range("C2").Value = 233
$ws.Range("C3").Value = 209
$ws.Range("C4").Value = 200
$ws.Range("C5").Value = 172
$ws.Range("C6").Value = 131
$ws.Range("C7").Value = 96

# Update Average (column D) values
$ws.Range("D2").Value = 266.25
$ws.Range("D3").Value = 222.5
$ws.Range("D4").Value = 241.5
$ws.Range("D5").Value = 233.25
$ws.Range("D6").Value = 119.4
$ws.Range("D7").Value = 80.8
